# Remove the trailing "Ver no Jupiter..." line and the copyright/footer
# line, along with the blank paragraph that separated them from the
# "Requisitos" section's "LOB1038: ..." requirement line. Everything else
# (including the blank paragraphs further down, before the page break)
# is left untouched.
$d = $word.ActiveDocument

$targets = @(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    [char]0x00A9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
)

# Walk from the end of the document backwards so that deleting a paragraph
# never invalidates the index of a paragraph we haven't visited yet.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]0x0007, [char]0x000D, [char]0x000A)

    foreach ($target in $targets) {
        if ($text -eq $target) {
            # Grab the (still valid, earlier-in-document) reference to the
            # preceding paragraph before mutating anything.
            $prev = $para.Previous()

            # Delete this paragraph first: it is the later range, so doing
            # this does not shift/invalidate $prev's range.
            $para.Range.Delete()

            if ($prev -ne $null) {
                $prevText = $prev.Range.Text.TrimEnd([char]0x0007, [char]0x000D, [char]0x000A)
                if ($prevText -eq "") {
                    $prev.Range.Delete()
                }
            }
            break
        }
    }
}
